$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate row 11 with the new "LocalRun" entry (values first, then copy the
# look of the existing bordered rows onto it so the new row matches the rest
# of the table).
$ws.Range("A11").Value = "LocalRun"
$ws.Range("B11").Value = "Running only API integration,ENW and LI scripts"
$ws.Range("C11").Value = "Y"

$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null

$ws.Range("B8").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null

$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Move / update the active selection like the original author's session.
$ws.Range("C15").Select() | Out-Null
